# PC adaption started, host does not work
#
# Moves the "8" entry in the "table" sheet from A8 down to A10,
# flips the card A2 marker from "+" to "-", and corrects the poster's
# event date in A10. Also nudges each sheet's active-cell selection to
# match where the author left off editing.

$wb = $excel.ActiveWorkbook

# --- "table" sheet ---------------------------------------------------
$wsTable = $wb.Worksheets.Item("table")
$wsTable.Range("A8").Value = $null
$wsTable.Range("A10").Value = 8
$wsTable.Range("E15").Select()

# --- "card" sheet ------------------------------------------------------
$wsCard = $wb.Worksheets.Item("card")
$wsCard.Range("A2").Value = "-"
$wsCard.Range("A2").Select()

# --- "poster" sheet ------------------------------------------------------
$wsPoster = $wb.Worksheets.Item("poster")
$wsPoster.Range("A10").Value = 45616
$wsPoster.Range("A11").Select()
